$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-02-15 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-16 Monday", 2)

$d.Content.Find.Execute("416÷2=208, 0", $true, $false, $false, $false, $false, $true, 1, $false, "821÷4=205, 1", 2)
$d.Content.Find.Execute("445÷3=148, 1", $true, $false, $false, $false, $false, $true, 1, $false, "225÷2=112, 1", 2)
$d.Content.Find.Execute("186÷5=37, 1", $true, $false, $false, $false, $false, $true, 1, $false, "252÷6=42, 0", 2)
$d.Content.Find.Execute("287÷2=143, 1", $true, $false, $false, $false, $false, $true, 1, $false, "900÷3=300, 0", 2)
$d.Content.Find.Execute("360÷3=120, 0", $true, $false, $false, $false, $false, $true, 1, $false, "963÷7=137, 4", 2)

$d.Content.Find.Execute("437÷3=145, 2", $true, $false, $false, $false, $false, $true, 1, $false, "318÷5=63, 3", 2)
$d.Content.Find.Execute("599÷2=299, 1", $true, $false, $false, $false, $false, $true, 1, $false, "342÷5=68, 2", 2)
$d.Content.Find.Execute("964÷4=241, 0", $true, $false, $false, $false, $false, $true, 1, $false, "255÷2=127, 1", 2)
$d.Content.Find.Execute("554÷8=69, 2", $true, $false, $false, $false, $false, $true, 1, $false, "378÷5=75, 3", 2)
$d.Content.Find.Execute("831÷3=277, 0", $true, $false, $false, $false, $false, $true, 1, $false, "725÷7=103, 4", 2)

$d.Content.Find.Execute("671÷2=335, 1", $true, $false, $false, $false, $false, $true, 1, $false, "956÷3=318, 2", 2)
$d.Content.Find.Execute("705÷4=176, 1", $true, $false, $false, $false, $false, $true, 1, $false, "460÷3=153, 1", 2)
$d.Content.Find.Execute("188÷2=94, 0", $true, $false, $false, $false, $false, $true, 1, $false, "604÷2=302, 0", 2)
$d.Content.Find.Execute("431÷3=143, 2", $true, $false, $false, $false, $false, $true, 1, $false, "983÷6=163, 5", 2)
$d.Content.Find.Execute("935÷7=133, 4", $true, $false, $false, $false, $false, $true, 1, $false, "499÷5=99, 4", 2)

$d.Content.Find.Execute("334÷8=41, 6", $true, $false, $false, $false, $false, $true, 1, $false, "972÷4=243, 0", 2)
$d.Content.Find.Execute("419÷4=104, 3", $true, $false, $false, $false, $false, $true, 1, $false, "544÷6=90, 4", 2)
$d.Content.Find.Execute("662÷9=73, 5", $true, $false, $false, $false, $false, $true, 1, $false, "871÷6=145, 1", 2)
$d.Content.Find.Execute("720÷3=240, 0", $true, $false, $false, $false, $false, $true, 1, $false, "805÷6=134, 1", 2)
$d.Content.Find.Execute("147÷6=24, 3", $true, $false, $false, $false, $false, $true, 1, $false, "162÷3=54, 0", 2)

$d.Content.Find.Execute("321÷3=107, 0", $true, $false, $false, $false, $false, $true, 1, $false, "549÷9=61, 0", 2)
$d.Content.Find.Execute("164÷3=54, 2", $true, $false, $false, $false, $false, $true, 1, $false, "878÷3=292, 2", 2)
$d.Content.Find.Execute("418÷3=139, 1", $true, $false, $false, $false, $false, $true, 1, $false, "900÷4=225, 0", 2)
$d.Content.Find.Execute("568÷9=63, 1", $true, $false, $false, $false, $false, $true, 1, $false, "624÷9=69, 3", 2)
$d.Content.Find.Execute("717÷4=179, 1", $true, $false, $false, $false, $false, $true, 1, $false, "160÷8=20, 0", 2)
